$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new TODO item in what is currently the next free row (row 5)
$ws.Cells.Item(5, 1).Value = "Add delay of ~2us between commutation or 30degrees"
$ws.Cells.Item(5, 2).Value = "to mitigate the noise, ringing of the voltage waveform"

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# New header cell in row 1
$ws.Cells.Item(1, 3).Value = "Status"

# New "DONE" status next to the first (now row 2) TODO item
$ws.Cells.Item(2, 3).Value = "DONE"

# Highlight the first TODO row (A2:B2) to show it's done
$ws.Range("A2:B2").Interior.ThemeColor = 9
$ws.Range("A2:B2").Interior.TintAndShade = 0.39997558519241921

# Adjust column A width to fit new content
$ws.Columns.Item(1).ColumnWidth = 50.5703125

# Update selection
$ws.Range("D2").Select()
